$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Trends Status"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("C2").Value = 0
$ws1.Range("E2").ClearContents()
$ws1.Range("E3").ClearContents()
$ws1.Range("E4").ClearContents()
$ws1.Range("E5").ClearContents()
$ws1.Range("E6").ClearContents()

$ws1.Range("B7").Value = 10
$ws1.Range("C7").Value = 48

$ws1.Range("B8").Value = 585
$ws1.Range("C8").Value = 547

# ---------------------------------------------------------------------------
# Sheet 3: "Priority Status"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Priority Status")

$ws3.Range("B2").Value = 103
$ws3.Range("B3").Value = 286
$ws3.Range("B4").Value = 554

# ---------------------------------------------------------------------------
# Sheet 4: "Species qualification"
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("A2").Value = "SoIB Assessment"
$ws4.Range("B2").Value = 595

$ws4.Range("B3").Value = 10

$ws4.Range("B4").Value = 48
$ws4.Range("C4").Value = 0

# ---------------------------------------------------------------------------
# Sheet 5: "High Priority break-up"
# First duplicate it (so the duplicate keeps the *old* numbers and becomes
# the new "Major update - High Priority " sheet), then rename the original
# and update its numbers to become "Interannual update - High Pri".
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("High Priority break-up")

$ws5.Copy([System.Reflection.Missing]::Value, $ws5)
$ws6 = $wb.Worksheets.Item($ws5.Index + 1)
$ws6.Name = "Major update - High Priority "

$ws5.Name = "Interannual update - High Pri"

$ws5.Range("B2").Value = 60
$ws5.Range("C2").Value = 58.3
$ws5.Range("D2").Value = 60
$ws5.Range("E2").Value = 82.2

$ws5.Range("B3").Value = 43
$ws5.Range("C3").Value = 41.7
$ws5.Range("D3").Value = 13
$ws5.Range("E3").Value = 17.8

# Restore the originally active/selected tab (sheet 1) so we don't introduce
# an unrelated sheetView change as a side effect of adding/copying sheets.
$ws1.Activate()
